$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.082.74"
$ws.Range("E2").Value = "  -1.38%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.637.15"
$ws.Range("E3").Value = "  -1.14%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9979"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.59"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9981"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3943"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3875"
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.35"
$ws.Range("E9").Value = "  +4.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.382"
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9985"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08523"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("E13").Value = "  -3.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.137"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("E16").Value = "  -0.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.615.76"
$ws.Range("E17").Value = "  -2.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.77"
$ws.Range("E18").Value = "  +1.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06916"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "20.21"
$ws.Range("E20").Value = "  -3.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.900"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9965"
$ws.Range("E22").Value = "  -0.46%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.53"
$ws.Range("E23").Value = "  -2.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.066.77"
$ws.Range("E24").Value = "  -1.41%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.456"
$ws.Range("E25").Value = "  +5.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.891"
$ws.Range("E26").Value = "  +3.93%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.37"
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "157.24"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "141.15"
$ws.Range("E29").Value = "  -2.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.356"
$ws.Range("E30").Value = "  -6.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.035"
$ws.Range("E31").Value = "  -1.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.530"
$ws.Range("E32").Value = "  -3.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.809.46"
$ws.Range("E33").Value = "  -1.59%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08180"
$ws.Range("E34").Value = "  -0.22%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.011"
$ws.Range("E35").Value = "  -0.85%  "
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.700"
$ws.Range("E37").Value = "  -2.49%  "
$ws.Range("E38").Value = "  -2.92%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "10.49"
$ws.Range("E39").Value = "  +2.42%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09208"
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("E41").Value = "  +2.86%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7621"
$ws.Range("E42").Value = "  -2.67%  "
$ws.Range("E43").Value = "  -3.76%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.26"
$ws.Range("E44").Value = "  -1.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6970"
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.489"
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.108"
$ws.Range("E47").Value = "  -1.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9972"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08351"
$ws.Range("E49").Value = "  -3.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.55"
$ws.Range("E50").Value = "  -0.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.452"
$ws.Range("E51").Value = "  +19.42%  "
